$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-15 (date shifted forward by 5 days, rolling window)
$data = @(
    @(45929, 5236.12347480818, 5315.33929068954, 7152, 6794.450112, -11.5972530049436),
    @(45930, 5220.96395511301, 5490.18063814705, 2952, 6775.748821, 170.540229334752),
    @(45931, 4177.55357137576, 5253.82814826742, 3764, 5977.141939, 137.059021495486),
    @(45932, 4219.81850078184, 5799.04791141544, 3764, 6054.860144, 161.253731443066),
    @(45933, 4345.23365733728, 5439.8627777128, 3764, 6294.674447, 151.054315307313),
    @(45934, $null, $null, 3764, 5588.458833, $null),
    @(45935, $null, $null, 3764, 5562.329802, $null),
    @(45936, $null, $null, 3764, 6192.581868, $null),
    @(45937, $null, $null, 3764, 6192.581868, $null),
    @(45938, $null, $null, 3764, 6192.581868, $null),
    @(45939, $null, $null, 3764, 6192.581868, $null),
    @(45940, $null, $null, 3764, 6192.581868, $null),
    @(45941, $null, $null, 3764, 5600.914911, $null),
    @(45942, $null, $null, 3764, 5574.78588, $null)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
